$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2021 data corrections / additions (column I) and new 2022 column (J) ---

$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 2

$ws.Range("J5").Value = 50000

$ws.Range("J7").Value = 650

$ws.Range("J12").Value = 1800

# Cobalt (row 14): previous mistake, we were taking refined cobalt and not raw cobalt mineral
$ws.Range("B14").Value = 3907
$ws.Range("C14").Value = 4339
$ws.Range("D14").Value = 4216
$ws.Range("E14").Value = 3704
$ws.Range("F14").Value = 3279
$ws.Range("G14").Value = 4365
$ws.Range("H14").Value = 4328
$ws.Range("I14").Value = 3964
$ws.Range("J14").Value = 3063
$ws.Range("K14").Value = "Previousmistake, we were taking refined cobalt and not raw cobalt mineral"

$ws.Range("I16").Value = 17600

$ws.Range("G19").Value = 90000
$ws.Range("H19").Value = 140000
$ws.Range("I19").Value = 140000
$ws.Range("J19").Value = 12000

$ws.Range("J21").Value = 205.831

$ws.Range("H22").Value = 7620
$ws.Range("I22").Value = 7700
$ws.Range("J22").Value = 13000

# Indium (row 25) comment
$ws.Range("K25").Value = "Approximated with refined indium production data"

$ws.Range("I28").Value = 34500000
$ws.Range("J28").Value = 41400000

$ws.Range("I31").Value = 10000
$ws.Range("J31").Value = 6000

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 515

$ws.Range("J36").Value = 15000

$ws.Range("J37").Value = 952

$ws.Range("H38").Value = 167243

$ws.Range("E39").Value = 7200
$ws.Range("F39").Value = 7600
$ws.Range("G39").Value = 6800
$ws.Range("H39").Value = 6400
$ws.Range("I39").Value = 7500

$ws.Range("G41").Value = 1327000
$ws.Range("H41").Value = 1379000
$ws.Range("I41").Value = 1653000
# Fix "Horticutural" typo -> "Horticultural"
$ws.Range("K41").Value = "Horticultural use"

$ws.Range("G45").Value = 12820000
$ws.Range("H45").Value = 14276000
$ws.Range("I45").Value = 14370000
$ws.Range("J45").Value = 13460000

$ws.Range("H53").Value = 13300000
$ws.Range("I53").Value = 11800000

$ws.Range("H57").Value = 153293
$ws.Range("J57").Value = 200000

$ws.Range("B66").Value = 3000
$ws.Range("C66").Value = 5600
$ws.Range("D66").Value = 12500

$ws.Range("I67").Value = 230000
$ws.Range("J67").Value = 208000

# Restore the view state (scroll position / active cell) recorded in the edit
$ws.Range("E54").Select()

Write-Output "applied minerals 2021/2022 data update"
